$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 430.46155
$ws.Range("I33").Value = 327.36365
$ws.Range("K33").Value = 327.36365
$ws.Range("M33").Value = -98.36365000000001
$ws.Range("H138").Value = 6707.8643
$ws.Range("J138").Value = 9026.380999999999
$ws.Range("L138").Value = 27079.143
$ws.Range("N138").Value = -37359.143

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3114.4443
$ws.Range("I61").Value = 2349.3333
$ws.Range("J61").Value = 4644.6665
$ws.Range("K61").Value = 2349.3333
$ws.Range("L61").Value = 4644.6665
$ws.Range("M61").Value = -2137.3333
$ws.Range("N61").Value = -5068.6665
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 2794.9412
$ws.Range("I132").Value = 2720.138
$ws.Range("J132").Value = 3228.8
$ws.Range("K132").Value = 8160.414
$ws.Range("L132").Value = 9686.400000000001
$ws.Range("M132").Value = -5630.414
$ws.Range("N132").Value = -14746.4
$ws.Range("H136").Value = 3114.4443
$ws.Range("I136").Value = 2349.3333
$ws.Range("J136").Value = 4644.6665
$ws.Range("K136").Value = 7047.999899999999
$ws.Range("L136").Value = 13933.9995
$ws.Range("M136").Value = -4497.999899999999
$ws.Range("N136").Value = -19033.9995

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3058.1875
$ws.Range("I105").Value = 2414
$ws.Range("J105").Value = 3886.4285
$ws.Range("K105").Value = 2414
$ws.Range("L105").Value = 3886.4285
$ws.Range("M105").Value = -667
$ws.Range("N105").Value = -7380.4285
$ws.Range("H134").Value = 1947.6842
$ws.Range("I134").Value = 1534.875
$ws.Range("K134").Value = 4604.625
$ws.Range("M134").Value = -2069.625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 588.6923
$ws.Range("I22").Value = 556.25
$ws.Range("J22").Value = 640.6
$ws.Range("K22").Value = 556.25
$ws.Range("L22").Value = 640.6
$ws.Range("M22").Value = -206.25
$ws.Range("N22").Value = -1340.6
$ws.Range("H31").Value = 6119.62
$ws.Range("I31").Value = 4003.9565
$ws.Range("J31").Value = 7921.852
$ws.Range("K31").Value = 4003.9565
$ws.Range("L31").Value = 7921.852
$ws.Range("M31").Value = -3708.9565
$ws.Range("N31").Value = -8511.851999999999
$ws.Range("H34").Value = 6119.62
$ws.Range("I34").Value = 4003.9565
$ws.Range("J34").Value = 7921.852
$ws.Range("K34").Value = 4003.9565
$ws.Range("L34").Value = 7921.852
$ws.Range("M34").Value = -3801.9565
$ws.Range("N34").Value = -8325.851999999999
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H60").Value = 8398.333000000001
$ws.Range("I60").Value = 97.5
$ws.Range("J60").Value = 25000
$ws.Range("K60").Value = 97.5
$ws.Range("L60").Value = 25000
$ws.Range("M60").Value = 413.5
$ws.Range("N60").Value = -26022
$ws.Range("H132").Value = 1433.75
$ws.Range("I132").Value = 1220.7
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 3662.1
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -1132.1
$ws.Range("N132").Value = -12557
$ws.Range("H138").Value = 78000
$ws.Range("J138").Value = 78000
$ws.Range("L138").Value = 78000
$ws.Range("N138").Value = -88280

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 416702
$ws.Range("J2").Value = 35.666668
$ws.Range("L2").Value = 214.000008
$ws.Range("N2").Value = -440.000008
$ws.Range("H7").Value = 589
$ws.Range("J7").Value = 198.5
$ws.Range("L7").Value = 595.5
$ws.Range("N7").Value = -819.5
$ws.Range("H12").Value = 1048.75
$ws.Range("I12").Value = 998.3333
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 2994.9999
$ws.Range("L12").Value = 3600
$ws.Range("M12").Value = -2821.9999
$ws.Range("N12").Value = -3946
$ws.Range("H55").Value = 875.25
$ws.Range("J55").Value = 999.3333
$ws.Range("L55").Value = 2997.9999
$ws.Range("N55").Value = -3351.9999
$ws.Range("H131").Value = 3839912
$ws.Range("I131").Value = 124343.22
$ws.Range("J131").Value = 6412229
$ws.Range("K131").Value = 373029.66
$ws.Range("L131").Value = 19236687
$ws.Range("M131").Value = -367989.66
$ws.Range("N131").Value = -19246767

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 593.8823
$ws.Range("J107").Value = 1019.8571
$ws.Range("L107").Value = 1019.8571
$ws.Range("N107").Value = -4859.8571
$ws.Range("H132").Value = 3131.05
$ws.Range("I132").Value = 3155.3572
$ws.Range("J132").Value = 3074.3333
$ws.Range("K132").Value = 9466.071599999999
$ws.Range("L132").Value = 9222.999899999999
$ws.Range("M132").Value = -6936.071599999999
$ws.Range("N132").Value = -14282.9999
$ws.Range("H135").Value = 189836.25
$ws.Range("I135").Value = 95000
$ws.Range("J135").Value = 284672.5
$ws.Range("K135").Value = 95000
$ws.Range("L135").Value = 284672.5
$ws.Range("M135").Value = -89930
$ws.Range("N135").Value = -294812.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 41257.25
$ws.Range("I38").Value = 30014.5
$ws.Range("J38").Value = 52500
$ws.Range("K38").Value = 30014.5
$ws.Range("L38").Value = 52500
$ws.Range("M38").Value = -29604.5
$ws.Range("N38").Value = -53320
$ws.Range("H132").Value = 2033.8462
$ws.Range("I132").Value = 1979.2
$ws.Range("J132").Value = 2216
$ws.Range("K132").Value = 5937.6
$ws.Range("L132").Value = 6648
$ws.Range("M132").Value = -3407.6
$ws.Range("N132").Value = -11708
$ws.Range("H136").Value = 1367.4546
$ws.Range("I136").Value = 780.5
$ws.Range("K136").Value = 2341.5
$ws.Range("M136").Value = 208.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 76248.5
$ws.Range("J46").Value = 76248.5
$ws.Range("L46").Value = 76248.5
$ws.Range("N46").Value = -76710.5
$ws.Range("H100").Value = 690.1429000000001
$ws.Range("I100").Value = 607.2
$ws.Range("K100").Value = 1214.4
$ws.Range("M100").Value = -673.4000000000001
$ws.Range("H101").Value = 10842.714
$ws.Range("J101").Value = 10842.714
$ws.Range("L101").Value = 10842.714
$ws.Range("N101").Value = -17332.714
$ws.Range("H122").Value = 7485.533
$ws.Range("I122").Value = 9261.817999999999
$ws.Range("K122").Value = 27785.454
$ws.Range("M122").Value = -25335.454
$ws.Range("H132").Value = 1788.2162
$ws.Range("J132").Value = 2532.8333
$ws.Range("L132").Value = 7598.499899999999
$ws.Range("N132").Value = -12658.4999
$ws.Range("H134").Value = 76248.5
$ws.Range("J134").Value = 76248.5
$ws.Range("L134").Value = 228745.5
$ws.Range("N134").Value = -233815.5
